$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Thu Jun 29 18:42:54 UTC 2023 with GitHub Actions
#
# Refresh the Price (D) and Volume(1h) (E) columns for each coin row.
# Rows 38/39 also swap places: the coin previously listed as
# "TrustWalletToken" now shows RenderToken's data and vice versa.
# The Price column (D) is forced to Text format before assignment so
# numeric-looking values (e.g. "1.0000", "233.71") are preserved exactly
# as text rather than being normalized into numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.568.01'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.855.06'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.71'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4738'
$ws.Range("E7").Value = '  +0.84%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2752'
$ws.Range("E8").Value = '  +1.74%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06326'
$ws.Range("E9").Value = '  -0.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '17.69'
$ws.Range("E10").Value = '  +8.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.886.92'
$ws.Range("E11").Value = '  +1.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07449'
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.002'
$ws.Range("E13").Value = '  +1.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '84.56'
$ws.Range("E14").Value = '  -0.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6263'
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.532.83'
$ws.Range("E16").Value = '  +1.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '244.42'
$ws.Range("E17").Value = '  +6.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.0000'
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.70'
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007334'
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.939'
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.938'
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.150'
$ws.Range("E24").Value = '  -0.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.92'
$ws.Range("E25").Value = '  -2.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.01'
$ws.Range("E26").Value = '  +1.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.879'
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("E28").Value = '  -1.42%  '
$ws.Range("E30").Value = '  -2.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.835'
$ws.Range("E31").Value = '  -1.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.04842'
$ws.Range("E32").Value = '  -0.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.137'
$ws.Range("E33").Value = '  -1.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7024'
$ws.Range("E34").Value = '  -0.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.692'
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.01895'
$ws.Range("E36").Value = '  +1.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.678'
$ws.Range("E37").Value = '  +1.73%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.000'
$ws.Range("E38").Value = '  +2.88%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8758'
$ws.Range("E39").Value = '  -3.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '106.73'
$ws.Range("E40").Value = '  +1.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9995'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.545'
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4060'
$ws.Range("E43").Value = '  -0.56%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.194'
$ws.Range("E44").Value = '  +1.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.74'
$ws.Range("E45").Value = '  +4.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1212'
$ws.Range("E46").Value = '  +1.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '33.59'
$ws.Range("E47").Value = '  +1.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.529'
$ws.Range("E48").Value = '  -0.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05537'
$ws.Range("E49").Value = '  -0.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.356'
$ws.Range("E50").Value = '  -2.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3678'
$ws.Range("E51").Value = '  +0.28%  '
